$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A55").Value = 'he; boyfriend'
$ws.Range("B55").Value = '彼|かれ'
$ws.Range("A56").Value = 'she; girlfriend'
$ws.Range("B56").Value = '彼女|かのじょ'
$ws.Range("A57").Value = 'they'
$ws.Range("B57").Value = '彼ら|かれら'
$ws.Range("A58").Value = 'boyfriend'
$ws.Range("B58").Value = '彼氏|かれし'
$ws.Range("A59").Value = 'age; era'
$ws.Range("B59").Value = '時代|じだい'
$ws.Range("A60").Value = 'electricity fee'
$ws.Range("B60").Value = '電気代|でんきだい'
$ws.Range("A61").Value = '90''s'
$ws.Range("B61").Value = '九十年代|きゅうじゅうねんだい'
$ws.Range("A62").Value = 'in one''s teens'
$ws.Range("B62").Value = '十代|じゅうだい'
$ws.Range("A63").Value = 'instead'
$ws.Range("B63").Value = '代わりに|かわりに'
$ws.Range("A64").Value = 'international students'
$ws.Range("B64").Value = '留学生|りゅうがくせい'
$ws.Range("A65").Value = 'to study abroad'
$ws.Range("B65").Value = '留学する|りゅうがくする'
$ws.Range("A66").Value = 'absence; not home'
$ws.Range("B66").Value = '留守|るす'
$ws.Range("A67").Value = 'family'
$ws.Range("B67").Value = '家族|かぞく'
$ws.Range("A68").Value = 'race; ethnic group'
$ws.Range("B68").Value = '民族|みんぞく'
$ws.Range("A69").Value = 'aquarium'
$ws.Range("B69").Value = '水族館|すいぞくかん'
$ws.Range("A70").Value = 'member of royalty'
$ws.Range("B70").Value = '王族|おうぞく'
$ws.Range("A71").Value = 'father'
$ws.Range("B71").Value = '父親|ちちおや'
$ws.Range("A72").Value = 'kind'
$ws.Range("B72").Value = '親切な|しんせつな'
$ws.Range("A73").Value = 'best friend'
$ws.Range("B73").Value = '親友|しんゆう'
$ws.Range("A74").Value = 'parents'
$ws.Range("B74").Value = '両親|りょうしん'
$ws.Range("A75").Value = 'intimate'
$ws.Range("B75").Value = '親しい|したしい'
$ws.Range("A76").Value = 'mother'
$ws.Range("B76").Value = '母親|ははおや'
$ws.Range("A77").Value = 'to cut'
$ws.Range("B77").Value = '切る|きる'
$ws.Range("A78").Value = 'ticket'
$ws.Range("B78").Value = '切符|きっぷ'
$ws.Range("A79").Value = 'postage stamp'
$ws.Range("B79").Value = '切手|きって'
$ws.Range("A80").Value = 'precious'
$ws.Range("B80").Value = '大切な|たいせつな'
$ws.Range("A81").Value = 'English language'
$ws.Range("B81").Value = '英語|えいご'
$ws.Range("A82").Value = 'United Kingdom'
$ws.Range("B82").Value = '英国|えいこく'
$ws.Range("A83").Value = 'English conversation'
$ws.Range("B83").Value = '英会話|えいかいわ'
$ws.Range("A84").Value = 'hero'
$ws.Range("B84").Value = '英雄|えいゆう'
$ws.Range("A85").Value = 'shop'
$ws.Range("B85").Value = '店|みせ'
$ws.Range("A86").Value = 'store clerk'
$ws.Range("B86").Value = '店員|てんいん'
$ws.Range("A87").Value = 'stall; kiosk'
$ws.Range("B87").Value = '売店|ばいてん'
$ws.Range("A88").Value = 'book store'
$ws.Range("B88").Value = '書店|しょてん'
$ws.Range("A89").Value = 'store manager'
$ws.Range("B89").Value = '店長|てんちょう'
$ws.Range("A90").Value = 'last year'
$ws.Range("B90").Value = '去年|きょねん'
$ws.Range("A91").Value = 'the past'
$ws.Range("B91").Value = '過去|かこ'
$ws.Range("A92").Value = 'to leave'
$ws.Range("B92").Value = '去る|さる'
$ws.Range("A93").Value = 'to erase'
$ws.Range("B93").Value = '消去する|しょうきょする'
$ws.Range("A94").Value = 'suddenly'
$ws.Range("B94").Value = '急に|きゅうに'
$ws.Range("A95").Value = 'to hurry'
$ws.Range("B95").Value = '急ぐ|いそぐ'
$ws.Range("A96").Value = 'express train'
$ws.Range("B96").Value = '急行|きゅうこう'
$ws.Range("A97").Value = 'super express'
$ws.Range("B97").Value = '特急|とっきゅう'
$ws.Range("A98").Value = 'to ride'
$ws.Range("B98").Value = '乗る|のる'
$ws.Range("A99").Value = 'vehicle'
$ws.Range("B99").Value = '乗り物|のりもの'
$ws.Range("A100").Value = 'riding a car'
$ws.Range("B100").Value = '乗車|じょうしゃ'
$ws.Range("A101").Value = 'horseback riding'
$ws.Range("B101").Value = '乗馬|じょうば'
$ws.Range("A102").Value = 'really'
$ws.Range("B102").Value = '本当に|ほんとうに'
$ws.Range("A103").Value = 'lunch box'
$ws.Range("B103").Value = 'お弁当|おべんとう'
$ws.Range("A104").Value = 'at that time'
$ws.Range("B104").Value = '当時|とうじ'
$ws.Range("A105").Value = 'to hit'
$ws.Range("B105").Value = '当たる|あたる'
$ws.Range("A106").Value = 'music'
$ws.Range("B106").Value = '音楽|おんがく'
$ws.Range("A107").Value = 'pronunciation'
$ws.Range("B107").Value = '発音|はつおん'
$ws.Range("A108").Value = 'sound'
$ws.Range("B108").Value = '音|おと'
$ws.Range("A109").Value = 'real intention'
$ws.Range("B109").Value = '本音|ほんね'
$ws.Range("A110").Value = 'fun'
$ws.Range("B110").Value = '楽しい|たのしい'
$ws.Range("A111").Value = 'musical instrument'
$ws.Range("B111").Value = '楽器|がっき'
$ws.Range("A112").Value = 'easy; comfortable'
$ws.Range("B112").Value = '楽な|らくな'
$ws.Range("A113").Value = 'doctor'
$ws.Range("B113").Value = '医者|いしゃ'
$ws.Range("A114").Value = 'dentist'
$ws.Range("B114").Value = '歯医者|はいしゃ'
$ws.Range("A115").Value = 'medical science'
$ws.Range("B115").Value = '医学|いがく'
$ws.Range("A116").Value = 'clinic'
$ws.Range("B116").Value = '医院|いいん'
$ws.Range("A117").Value = 'scholar'
$ws.Range("B117").Value = '学者|がくしゃ'
$ws.Range("A118").Value = 'reader'
$ws.Range("B118").Value = '読者|どくしゃ'
$ws.Range("A119").Value = 'young people'
$ws.Range("B119").Value = '若者|わかもの'
$ws.Range("A120").Value = 'ninja'
$ws.Range("B120").Value = '忍者|にんじゃ'
